$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: value 0, bold font, thin border all around, centered horizontally, top-aligned vertically
$b1 = $ws.Range("B1")
$b1.Value = 0
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160
$b1.Borders.LineStyle = 1
$b1.Borders.Weight = 2

# Copy the same formatting onto A2 (reuses the style just created instead of
# re-deriving it incrementally, which would otherwise leave a stray unused xf)
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)
$a2.Value = 0

# B2: plain text value (becomes a shared string), default formatting
$ws.Range("B2").Value = "disconnected_elements"
